# Weekly price update: a new week's record is inserted as row 17 (just
# below the most-recent-so-far block), pushing the previously-existing
# rows 17:40 down to 18:41. The sheet's used range grows from A1:R40 to
# A1:R41 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 17-40 down to 18-41, carrying formatting
# (incl. the date number format on column D) down with them.
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with this week's record.
$ws.Cells.Item(17, 1).Value = 8
$ws.Cells.Item(17, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(17, 3).Value = "Coquimbo"
$ws.Cells.Item(17, 4).Value = 44880
$ws.Cells.Item(17, 5).Value = 4
$ws.Cells.Item(17, 6).Value = 100112026
$ws.Cells.Item(17, 7).Value = "Haba"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 560
$ws.Cells.Item(17, 11).Value = 7000
$ws.Cells.Item(17, 12).Value = 8000
$ws.Cells.Item(17, 13).Value = 7500
$ws.Cells.Item(17, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(17, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(17, 16).Value = 300
$ws.Cells.Item(17, 17).Value = 25
$ws.Cells.Item(17, 18).Value = "Hortaliza"
